$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected by Excel as a number
# (pure digits + one decimal point) are first marked as Text format so the
# assignment keeps them as strings (matching the source inlineStr cells),
# then the style is reset back to Normal so no stray formatting is left behind.
$forceTextCells = @("D5", "D6", "D14", "D20", "D21", "D22", "D25", "D27", "D31", "D34", "D38", "D39", "D42", "D43", "D47", "D48")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.642.78"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.703.44"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "599.28"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "160.40"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "2.702.68"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "28.33"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "3.196.20"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "68.572.98"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "2.703.27"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("D20").Value = "366.27"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "7.63"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").Value = "74.52"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "9.96"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("E30").Value = "  -6.48%  "
$ws.Range("D31").Value = "578.98"
$ws.Range("E31").Value = "  +4.51%  "
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("D34").Value = "1.95"
$ws.Range("E34").Value = "  +6.10%  "
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +6.50%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "19.88"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").Value = "160.67"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").Value = "5.41"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "2.71"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "0.0₆0319"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("D47").Value = "158.00"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "3.96"
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("E50").Value = "  +7.02%  "
$ws.Range("E51").Value = "  -0.12%  "

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
